$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (2-215).
for ($r = 2; $r -le 215; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value2 = 45175
    }
}
